$d = $word.ActiveDocument
$d.Content.Find.Execute("11+54=65", $true, $false, $false, $false, $false, $false, 1, $false, "48-24=24", 2)
$d.Content.Find.Execute("39-39=0", $true, $false, $false, $false, $false, $false, 1, $false, "8+38=46", 2)
$d.Content.Find.Execute("53+14=67", $true, $false, $false, $false, $false, $false, 1, $false, "82-13=69", 2)
$d.Content.Find.Execute("96-74=22", $true, $false, $false, $false, $false, $false, 1, $false, "75-62=13", 2)
$d.Content.Find.Execute("34+30=64", $true, $false, $false, $false, $false, $false, 1, $false, "26-7=19", 2)
$d.Content.Find.Execute("22+13=35", $true, $false, $false, $false, $false, $false, 1, $false, "27+48=75", 2)
$d.Content.Find.Execute("24+15=39", $true, $false, $false, $false, $false, $false, 1, $false, "61+33=94", 2)
$d.Content.Find.Execute("18-6=12", $true, $false, $false, $false, $false, $false, 1, $false, "10+57=67", 2)
$d.Content.Find.Execute("48-8=40", $true, $false, $false, $false, $false, $false, 1, $false, "84-6=78", 2)
$d.Content.Find.Execute("28+6=34", $true, $false, $false, $false, $false, $false, 1, $false, "12+82=94", 2)
$d.Content.Find.Execute("53-31=22", $true, $false, $false, $false, $false, $false, 1, $false, "85-39=46", 2)
$d.Content.Find.Execute("30-22=8", $true, $false, $false, $false, $false, $false, 1, $false, "29+62=91", 2)
$d.Content.Find.Execute("79+1=80", $true, $false, $false, $false, $false, $false, 1, $false, "83-72=11", 2)
$d.Content.Find.Execute("30+24=54", $true, $false, $false, $false, $false, $false, 1, $false, "72-42=30", 2)
$d.Content.Find.Execute("87-41=46", $true, $false, $false, $false, $false, $false, 1, $false, "59-38=21", 2)
$d.Content.Find.Execute("40+28=68", $true, $false, $false, $false, $false, $false, 1, $false, "70-54=16", 2)
$d.Content.Find.Execute("88-37=51", $true, $false, $false, $false, $false, $false, 1, $false, "37-1=36", 2)
$d.Content.Find.Execute("60-54=6", $true, $false, $false, $false, $false, $false, 1, $false, "60-15=45", 2)
$d.Content.Find.Execute("45+1=46", $true, $false, $false, $false, $false, $false, 1, $false, "57-44=13", 2)
$d.Content.Find.Execute("99-20=79", $true, $false, $false, $false, $false, $false, 1, $false, "55-6=49", 2)
$d.Content.Find.Execute("52+28=80", $true, $false, $false, $false, $false, $false, 1, $false, "65-60=5", 2)
$d.Content.Find.Execute("21+52=73", $true, $false, $false, $false, $false, $false, 1, $false, "68+6=74", 2)
$d.Content.Find.Execute("21+74=95", $true, $false, $false, $false, $false, $false, 1, $false, "36+62=98", 2)
$d.Content.Find.Execute("28+50=78", $true, $false, $false, $false, $false, $false, 1, $false, "89-20=69", 2)
$d.Content.Find.Execute("28-11=17", $true, $false, $false, $false, $false, $false, 1, $false, "92-86=6", 2)
$d.Content.Find.Execute("58-47=11", $true, $false, $false, $false, $false, $false, 1, $false, "20+65=85", 2)
$d.Content.Find.Execute("54-52=2", $true, $false, $false, $false, $false, $false, 1, $false, "78-40=38", 2)
$d.Content.Find.Execute("16+37=53", $true, $false, $false, $false, $false, $false, 1, $false, "23-12=11", 2)
$d.Content.Find.Execute("61-22=39", $true, $false, $false, $false, $false, $false, 1, $false, "37-36=1", 2)
$d.Content.Find.Execute("29+33=62", $true, $false, $false, $false, $false, $false, 1, $false, "9+23=32", 2)
$d.Content.Find.Execute("55-15=40", $true, $false, $false, $false, $false, $false, 1, $false, "7+3=10", 2)
$d.Content.Find.Execute("13+66=79", $true, $false, $false, $false, $false, $false, 1, $false, "50-15=35", 2)
$d.Content.Find.Execute("14+72=86", $true, $false, $false, $false, $false, $false, 1, $false, "55-14=41", 2)
$d.Content.Find.Execute("54-27=27", $true, $false, $false, $false, $false, $false, 1, $false, "52-34=18", 2)
$d.Content.Find.Execute("35+21=56", $true, $false, $false, $false, $false, $false, 1, $false, "22+1=23", 2)
$d.Content.Find.Execute("75-25=50", $true, $false, $false, $false, $false, $false, 1, $false, "76-50=26", 2)
$d.Content.Find.Execute("1+40=41", $true, $false, $false, $false, $false, $false, 1, $false, "20-19=1", 2)
$d.Content.Find.Execute("61-21=40", $true, $false, $false, $false, $false, $false, 1, $false, "80+0=80", 2)
$d.Content.Find.Execute("69-8=61", $true, $false, $false, $false, $false, $false, 1, $false, "83-48=35", 2)
$d.Content.Find.Execute("0+22=22", $true, $false, $false, $false, $false, $false, 1, $false, "28+9=37", 2)
$d.Content.Find.Execute("59+24=83", $true, $false, $false, $false, $false, $false, 1, $false, "10-4=6", 2)
$d.Content.Find.Execute("27-7=20", $true, $false, $false, $false, $false, $false, 1, $false, "64-37=27", 2)
$d.Content.Find.Execute("90-42=48", $true, $false, $false, $false, $false, $false, 1, $false, "60-18=42", 2)
$d.Content.Find.Execute("55+11=66", $true, $false, $false, $false, $false, $false, 1, $false, "10+79=89", 2)
$d.Content.Find.Execute("74-48=26", $true, $false, $false, $false, $false, $false, 1, $false, "0+31=31", 2)
$d.Content.Find.Execute("43-41=2", $true, $false, $false, $false, $false, $false, 1, $false, "90-69=21", 2)
$d.Content.Find.Execute("9+10=19", $true, $false, $false, $false, $false, $false, 1, $false, "73-64=9", 2)
$d.Content.Find.Execute("42+33=75", $true, $false, $false, $false, $false, $false, 1, $false, "62-22=40", 2)
$d.Content.Find.Execute("81+9=90", $true, $false, $false, $false, $false, $false, 1, $false, "28+8=36", 2)
$d.Content.Find.Execute("62-49=13", $true, $false, $false, $false, $false, $false, 1, $false, "11+68=79", 2)
$d.Content.Find.Execute("29+65=94", $true, $false, $false, $false, $false, $false, 1, $false, "70-50=20", 2)
$d.Content.Find.Execute("45-3=42", $true, $false, $false, $false, $false, $false, 1, $false, "37+3=40", 2)
$d.Content.Find.Execute("10+82=92", $true, $false, $false, $false, $false, $false, 1, $false, "42+19=61", 2)
$d.Content.Find.Execute("93-63=30", $true, $false, $false, $false, $false, $false, 1, $false, "25-11=14", 2)
$d.Content.Find.Execute("48+31=79", $true, $false, $false, $false, $false, $false, 1, $false, "96-31=65", 2)
$d.Content.Find.Execute("71-26=45", $true, $false, $false, $false, $false, $false, 1, $false, "47+9=56", 2)
$d.Content.Find.Execute("24+27=51", $true, $false, $false, $false, $false, $false, 1, $false, "26+6=32", 2)
$d.Content.Find.Execute("0+1=1", $true, $false, $false, $false, $false, $false, 1, $false, "24-20=4", 2)
$d.Content.Find.Execute("49+23=72", $true, $false, $false, $false, $false, $false, 1, $false, "41-27=14", 2)
$d.Content.Find.Execute("51+24=75", $true, $false, $false, $false, $false, $false, 1, $false, "77-4=73", 2)
$d.Content.Find.Execute("35-15=20", $true, $false, $false, $false, $false, $false, 1, $false, "17+76=93", 2)
$d.Content.Find.Execute("1+52=53", $true, $false, $false, $false, $false, $false, 1, $false, "54+45=99", 2)
$d.Content.Find.Execute("69-39=30", $true, $false, $false, $false, $false, $false, 1, $false, "19+53=72", 2)
$d.Content.Find.Execute("20+55=75", $true, $false, $false, $false, $false, $false, 1, $false, "8+75=83", 2)
$d.Content.Find.Execute("20+43=63", $true, $false, $false, $false, $false, $false, 1, $false, "92-39=53", 2)
$d.Content.Find.Execute("47+23=70", $true, $false, $false, $false, $false, $false, 1, $false, "76-40=36", 2)
$d.Content.Find.Execute("99-66=33", $true, $false, $false, $false, $false, $false, 1, $false, "53-47=6", 2)
$d.Content.Find.Execute("14+85=99", $true, $false, $false, $false, $false, $false, 1, $false, "43+34=77", 2)
$d.Content.Find.Execute("47-28=19", $true, $false, $false, $false, $false, $false, 1, $false, "57+19=76", 2)
$d.Content.Find.Execute("42+25=67", $true, $false, $false, $false, $false, $false, 1, $false, "59-23=36", 2)
$d.Content.Find.Execute("20+18=38", $true, $false, $false, $false, $false, $false, 1, $false, "70-21=49", 2)
$d.Content.Find.Execute("94-90=4", $true, $false, $false, $false, $false, $false, 1, $false, "93-14=79", 2)
$d.Content.Find.Execute("71+15=86", $true, $false, $false, $false, $false, $false, 1, $false, "30+6=36", 2)
$d.Content.Find.Execute("48-25=23", $true, $false, $false, $false, $false, $false, 1, $false, "10+35=45", 2)
$d.Content.Find.Execute("46-30=16", $true, $false, $false, $false, $false, $false, 1, $false, "26+73=99", 2)
$d.Content.Find.Execute("4+74=78", $true, $false, $false, $false, $false, $false, 1, $false, "18+15=33", 2)
$d.Content.Find.Execute("31+4=35", $true, $false, $false, $false, $false, $false, 1, $false, "66-22=44", 2)
$d.Content.Find.Execute("75+15=90", $true, $false, $false, $false, $false, $false, 1, $false, "16+41=57", 2)
$d.Content.Find.Execute("60-36=24", $true, $false, $false, $false, $false, $false, 1, $false, "8+62=70", 2)
$d.Content.Find.Execute("58+27=85", $true, $false, $false, $false, $false, $false, 1, $false, "6+3=9", 2)
$d.Content.Find.Execute("48+22=70", $true, $false, $false, $false, $false, $false, 1, $false, "27+9=36", 2)
$d.Content.Find.Execute("6+24=30", $true, $false, $false, $false, $false, $false, 1, $false, "27-21=6", 2)
$d.Content.Find.Execute("23+41=64", $true, $false, $false, $false, $false, $false, 1, $false, "55-42=13", 2)
$d.Content.Find.Execute("82-3=79", $true, $false, $false, $false, $false, $false, 1, $false, "14-7=7", 2)
$d.Content.Find.Execute("97-19=78", $true, $false, $false, $false, $false, $false, 1, $false, "94-32=62", 2)
$d.Content.Find.Execute("98-28=70", $true, $false, $false, $false, $false, $false, 1, $false, "91-5=86", 2)
$d.Content.Find.Execute("37+58=95", $true, $false, $false, $false, $false, $false, 1, $false, "6+89=95", 2)
$d.Content.Find.Execute("37-34=3", $true, $false, $false, $false, $false, $false, 1, $false, "80-54=26", 2)
$d.Content.Find.Execute("47-15=32", $true, $false, $false, $false, $false, $false, 1, $false, "20-0=20", 2)
$d.Content.Find.Execute("97-23=74", $true, $false, $false, $false, $false, $false, 1, $false, "93-24=69", 2)
$d.Content.Find.Execute("1+32=33", $true, $false, $false, $false, $false, $false, 1, $false, "64-8=56", 2)
$d.Content.Find.Execute("47+10=57", $true, $false, $false, $false, $false, $false, 1, $false, "0+14=14", 2)
$d.Content.Find.Execute("17+35=52", $true, $false, $false, $false, $false, $false, 1, $false, "36+33=69", 2)
$d.Content.Find.Execute("51-25=26", $true, $false, $false, $false, $false, $false, 1, $false, "34+25=59", 2)
$d.Content.Find.Execute("21+42=63", $true, $false, $false, $false, $false, $false, 1, $false, "4+5=9", 2)
$d.Content.Find.Execute("47+41=88", $true, $false, $false, $false, $false, $false, 1, $false, "25+19=44", 2)
$d.Content.Find.Execute("48-26=22", $true, $false, $false, $false, $false, $false, 1, $false, "0+51=51", 2)
$d.Content.Find.Execute("39-16=23", $true, $false, $false, $false, $false, $false, 1, $false, "2+91=93", 2)
$d.Content.Find.Execute("11+21=32", $true, $false, $false, $false, $false, $false, 1, $false, "54+22=76", 2)
$d.Content.Find.Execute("74-22=52", $true, $false, $false, $false, $false, $false, 1, $false, "59-4=55", 2)
Write-Output "done"
